$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 1.55
$ws.Range("K2").Value = 4.6
$ws.Range("U2").Value = 1.9
$ws.Range("V2").Value = 2.48
$ws.Range("F3").Value = 4.7
$ws.Range("K3").Value = 3.9
$ws.Range("K4").Value = 5
$ws.Range("AH4").Value = 50
$ws.Range("Q5").Value = 1.6
$ws.Range("N10").Value = 4.7
$ws.Range("P10").Value = 2.04
$ws.Range("G11").Value = 21
$ws.Range("F13").Value = 4.4
$ws.Range("G13").Value = 4.5
$ws.Range("H13").Value = 1.84
$ws.Range("I13").Value = 1.86
$ws.Range("K13").Value = 4.3
$ws.Range("AB13").Value = 18
$ws.Range("AJ13").Value = 100
$ws.Range("F14").Value = 1.32
$ws.Range("G14").Value = 1.41
$ws.Range("H14").Value = 9
$ws.Range("I14").Value = 11.5
$ws.Range("J14").Value = 5.4
$ws.Range("K14").Value = 6.6
$ws.Range("N14").Value = 2.64
$ws.Range("O14").Value = 1.16
$ws.Range("P14").Value = 2.64
$ws.Range("R14").Value = 1.57
$ws.Range("S14").Value = 2.04
$ws.Range("T14").Value = 1.66
$ws.Range("U14").Value = 1.7
$ws.Range("V14").Value = 1.09
$ws.Range("W14").Value = 3.4
$ws.Range("AB14").Value = 17
$ws.Range("AF14").Value = 12.5
$ws.Range("AG14").Value = 13.5
$ws.Range("AJ14").Value = 15
$ws.Range("AK14").Value = 18
$ws.Range("H15").Value = 10.5
$ws.Range("U15").Value = 1.8
$ws.Range("Y15").Value = 30
$ws.Range("AA15").Value = 500
$ws.Range("AB15").Value = 8
$ws.Range("AC15").Value = 12
$ws.Range("AE15").Value = 190
$ws.Range("AI15").Value = 160
$ws.Range("AK15").Value = 14.5
$ws.Range("J16").Value = 3.5
$ws.Range("K16").Value = 3.55
$ws.Range("P16").Value = 1.8
$ws.Range("Q16").Value = 2.2
$ws.Range("U16").Value = 2.02
$ws.Range("AO16").Value = 23
$ws.Range("L17").Value = 1.32
$ws.Range("T17").Value = 1.77
$ws.Range("H18").Value = 1.87
$ws.Range("I18").Value = 1.88
$ws.Range("J18").Value = 3.85
$ws.Range("K18").Value = 3.9
$ws.Range("P18").Value = 1.87
$ws.Range("T18").Value = 1.97
$ws.Range("V18").Value = 2.12
$ws.Range("AA18").Value = 20
$ws.Range("AC18").Value = 8.199999999999999
$ws.Range("AO18").Value = 14
$ws.Range("G19").Value = 1.73
$ws.Range("N19").Value = 3.55
$ws.Range("P19").Value = 2.02
$ws.Range("Q19").Value = 1.73
$ws.Range("S19").Value = 2.22
$ws.Range("K20").Value = 3.75
$ws.Range("L20").Value = 1.37
$ws.Range("M20").Value = 1.05
$ws.Range("N20").Value = 4
$ws.Range("O20").Value = 1.28
$ws.Range("R20").Value = 1.35
$ws.Range("S20").Value = 2.8
$ws.Range("T20").Value = 1.68
$ws.Range("U20").Value = 2.28
$ws.Range("V20").Value = 1.56
$ws.Range("W20").Value = 1.53
$ws.Range("X20").Value = 18.5
$ws.Range("Y20").Value = 12.5
$ws.Range("Z20").Value = 19.5
$ws.Range("AA20").Value = 42
$ws.Range("AB20").Value = 14.5
$ws.Range("AC20").Value = 9.4
$ws.Range("AD20").Value = 14.5
$ws.Range("AE20").Value = 29
$ws.Range("AF20").Value = 19.5
$ws.Range("AG20").Value = 13
$ws.Range("AH20").Value = 16
$ws.Range("AI20").Value = 38
$ws.Range("AJ20").Value = 42
$ws.Range("AK20").Value = 29
$ws.Range("AL20").Value = 38
$ws.Range("AM20").Value = 80
$ws.Range("AN20").Value = 23
$ws.Range("AO20").Value = 23
$ws.Range("F21").Value = 4.8
$ws.Range("G21").Value = 4.9
$ws.Range("I21").Value = 1.81
$ws.Range("K21").Value = 4.7
$ws.Range("L21").Value = 1.34
$ws.Range("Q21").Value = 1.76
$ws.Range("T21").Value = 1.79
$ws.Range("V21").Value = 2.24
$ws.Range("W21").Value = 1.25
$ws.Range("AB21").Value = 21
$ws.Range("AG21").Value = 23
$ws.Range("AI21").Value = 40
$ws.Range("AK21").Value = 75
$ws.Range("AL21").Value = 75
$ws.Range("AN21").Value = 85
$ws.Range("O22").Value = 1.26
$ws.Range("Q22").Value = 1.78
$ws.Range("R22").Value = 1.43
$ws.Range("S22").Value = 2.96
$ws.Range("T22").Value = 2.1
$ws.Range("U22").Value = 1.78
$ws.Range("Y22").Value = 1000
$ws.Range("AC22").Value = 12.5
$ws.Range("K23").Value = 15.5
$ws.Range("R23").Value = 2.18
$ws.Range("S23").Value = 1.75
$ws.Range("U23").Value = 1.67
$ws.Range("AL23").Value = 55
$ws.Range("G24").Value = 2.38
$ws.Range("H24").Value = 3.35
$ws.Range("I24").Value = 3.7
$ws.Range("J24").Value = 3.4
$ws.Range("L24").Value = 1.42
$ws.Range("O24").Value = 1.32
$ws.Range("P24").Value = 1.92
$ws.Range("Q24").Value = 1.95
$ws.Range("V24").Value = 1.37
$ws.Range("W24").Value = 1.72
$ws.Range("F25").Value = 1.39
$ws.Range("G25").Value = 1.4
$ws.Range("H25").Value = 9
$ws.Range("K25").Value = 5.8
$ws.Range("L25").Value = 1.28
$ws.Range("T25").Value = 1.92
$ws.Range("U25").Value = 2.02
$ws.Range("V25").Value = 1.11
$ws.Range("W25").Value = 3.5
$ws.Range("X25").Value = 24
$ws.Range("AC25").Value = 12.5
$ws.Range("AF25").Value = 8.800000000000001
$ws.Range("AH25").Value = 24
$ws.Range("AI25").Value = 110
$ws.Range("AL25").Value = 30
$ws.Range("AM25").Value = 120
$ws.Range("AO25").Value = 120
$ws.Range("L26").Value = 1.39
$ws.Range("N26").Value = 4.2
$ws.Range("U26").Value = 2.14
$ws.Range("X26").Value = 16.5
$ws.Range("AA26").Value = 19
$ws.Range("AC26").Value = 8.800000000000001
$ws.Range("AD26").Value = 9.800000000000001
$ws.Range("AE26").Value = 18
$ws.Range("AH26").Value = 19
$ws.Range("AI26").Value = 32
$ws.Range("AJ26").Value = 110
$ws.Range("AK26").Value = 55
$ws.Range("AM26").Value = 95
$ws.Range("AN26").Value = 60
$ws.Range("AO26").Value = 11
$ws.Range("F28").Value = 1.96
$ws.Range("G28").Value = 1.98
$ws.Range("H28").Value = 4.3
$ws.Range("F29").Value = 2.2
$ws.Range("J29").Value = 3.35
$ws.Range("K29").Value = 3.65
$ws.Range("F30").Value = 4.9
$ws.Range("H30").Value = 1.78
$ws.Range("K30").Value = 3.95
$ws.Range("K31").Value = 8.199999999999999
$ws.Range("F32").Value = 2.6
$ws.Range("I32").Value = 3.3
$ws.Range("J32").Value = 3.15
$ws.Range("K32").Value = 3.3
$ws.Range("P32").Value = 1.61
$ws.Range("Q32").Value = 2.44
